$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff: refreshed Price (D) / Volume(1h) (E) values,
# and for rows 17-24 the Coin (B) / Link (C) columns also rotate as the source list re-sorted.
$updates = @(
    @{ Ref = "D2"; Value = "331.55" },
    @{ Ref = "E2"; Value = "0.69%" },
    @{ Ref = "D3"; Value = "45.43" },
    @{ Ref = "E3"; Value = "3.19%" },
    @{ Ref = "D4"; Value = "5.608" },
    @{ Ref = "E4"; Value = "2.19%" },
    @{ Ref = "D5"; Value = "0.08357" },
    @{ Ref = "E5"; Value = "4.66%" },
    @{ Ref = "D6"; Value = "2.064" },
    @{ Ref = "E6"; Value = "4.15%" },
    @{ Ref = "D7"; Value = "0.9625" },
    @{ Ref = "E7"; Value = "1.33%" },
    @{ Ref = "D8"; Value = "2.583" },
    @{ Ref = "E8"; Value = "0.34%" },
    @{ Ref = "E9"; Value = "5.33%" },
    @{ Ref = "D10"; Value = "0.1923" },
    @{ Ref = "E10"; Value = "2.02%" },
    @{ Ref = "D11"; Value = "10.44" },
    @{ Ref = "E11"; Value = "-1.64%" },
    @{ Ref = "D12"; Value = "0.09883" },
    @{ Ref = "E12"; Value = "-0.94%" },
    @{ Ref = "D13"; Value = "0.04617" },
    @{ Ref = "E13"; Value = "-4.35%" },
    @{ Ref = "D14"; Value = "0.1060" },
    @{ Ref = "E14"; Value = "-0.26%" },
    @{ Ref = "D15"; Value = "0.001285" },
    @{ Ref = "E15"; Value = "0.41%" },
    @{ Ref = "D16"; Value = "0.006083" },
    @{ Ref = "E16"; Value = "2.47%" },
    @{ Ref = "B17"; Value = "LEO" },
    @{ Ref = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Ref = "D17"; Value = "3.377" },
    @{ Ref = "E17"; Value = "0.36%" },
    @{ Ref = "B18"; Value = "GateToken" },
    @{ Ref = "C18"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Ref = "D18"; Value = "4.441" },
    @{ Ref = "E18"; Value = "1.52%" },
    @{ Ref = "B19"; Value = "BitpandaEcosystemToken" },
    @{ Ref = "C19"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" },
    @{ Ref = "D19"; Value = "0.3369" },
    @{ Ref = "E19"; Value = "-2.87%" },
    @{ Ref = "B20"; Value = "ProBitToken" },
    @{ Ref = "C20"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" },
    @{ Ref = "D20"; Value = "0.1393" },
    @{ Ref = "E20"; Value = "-1.93%" },
    @{ Ref = "B21"; Value = "ZBToken" },
    @{ Ref = "C21"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb" },
    @{ Ref = "D21"; Value = "0.2654" },
    @{ Ref = "E21"; Value = "2.55%" },
    @{ Ref = "B22"; Value = "CoinExToken" },
    @{ Ref = "C22"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" },
    @{ Ref = "D22"; Value = "0.04187" },
    @{ Ref = "E22"; Value = "2.42%" },
    @{ Ref = "B23"; Value = "BitKan" },
    @{ Ref = "C23"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan" },
    @{ Ref = "D23"; Value = "0.001315" },
    @{ Ref = "E23"; Value = "3.67%" },
    @{ Ref = "B24"; Value = "HotbitToken" },
    @{ Ref = "C24"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" },
    @{ Ref = "D24"; Value = "0.004566" },
    @{ Ref = "E24"; Value = "7.04%" },
    @{ Ref = "E25"; Value = "8.56%" },
    @{ Ref = "D26"; Value = "0.0003750" },
    @{ Ref = "E26"; Value = "0.07%" },
    @{ Ref = "D38"; Value = "0.02714" },
    @{ Ref = "E38"; Value = "5.46%" },
    @{ Ref = "D39"; Value = "0.05766" },
    @{ Ref = "E39"; Value = "2.17%" },
    @{ Ref = "D40"; Value = "0.007841" },
    @{ Ref = "E40"; Value = "3.78%" },
    @{ Ref = "D41"; Value = "0.1434" },
    @{ Ref = "E41"; Value = "2.64%" },
    @{ Ref = "D42"; Value = "0.007294" },
    @{ Ref = "E42"; Value = "-1.75%" },
    @{ Ref = "D43"; Value = "0.002017" },
    @{ Ref = "E43"; Value = "-0.03%" },
    @{ Ref = "D44"; Value = "0.009090" },
    @{ Ref = "E44"; Value = "5.60%" },
    @{ Ref = "D45"; Value = "0.3549" },
    @{ Ref = "D46"; Value = "0.00007153" },
    @{ Ref = "E46"; Value = "0.17%" },
    @{ Ref = "E47"; Value = "0.19%" },
    @{ Ref = "D48"; Value = "0.0005816" },
    @{ Ref = "E48"; Value = "0.08%" },
    @{ Ref = "D49"; Value = "0.003501" },
    @{ Ref = "E49"; Value = "-7.29%" },
    @{ Ref = "D50"; Value = "0.003508" },
    @{ Ref = "E50"; Value = "-0.72%" },
    @{ Ref = "E51"; Value = "0.19%" }
)

foreach ($u in $updates) {
    $col = $u.Ref.Substring(0, 1)
    if ($col -eq "D" -or $col -eq "E") {
        # Force text format so numeric-looking / percent-looking strings are not
        # re-interpreted as numbers, matching the inline-string text cells in the workbook.
        $ws.Range($u.Ref).NumberFormat = "@"
    }
    $ws.Range($u.Ref).Value = $u.Value
}
